$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" column header in S4, reusing R4's style (s=21).
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# Q5/S5 pick up R5's existing style (s=22); R5 keeps it and gets a new value too.
$ws.Range("R5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("S5").PasteSpecial(-4122)

$ws.Range("Q5").Value = 91.892815141492093
$ws.Range("R5").Value = 101.53074848578628
$ws.Range("S5").Value = 109.27053140096621

# Update the active selection to match the authored sheetView.
$ws.Range("T5").Select()
